# Form the consolidated report: populate the "Absent" column (H) values
# that were left blank/zero, deriving Absent = 1 - Real (column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 6, 9, 14, 15, 16, 19, 20)

foreach ($r in $rows) {
    $real = $ws.Cells.Item($r, 5).Value2   # column E: Real
    if ($real -eq 1) {
        $ws.Cells.Item($r, 8).Value = 0
    } else {
        $ws.Cells.Item($r, 8).Value = 1
    }
}
